# Update cryptocurrency price/volume data on the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" column (D) values are stored as text (they use "." as both thousands
# separator and decimal point, which is not a valid Excel number), so force the
# NumberFormat to Text before assigning to prevent Excel from re-interpreting them
# as numbers.
$priceCells = @("D2", "D3", "D6", "D8", "D9", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D31", "D33", "D35", "D36", "D38", "D40", "D42", "D46", "D47", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "68.780.46"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.436.80"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "161.58"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.512"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "0.168"
$ws.Range("E9").Value = "  +6.95%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "4.58"
$ws.Range("E12").Value = "  -6.00%  "
$ws.Range("D13").Value = "0.0000177"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "68.692.92"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "2.884.32"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "23.22"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "2.435.21"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "10.54"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "339.01"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "6.94"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "3.84"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "1.93"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "66.94"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "3.71"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "2.561.42"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "8.22"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").Value = "0.0₃0819"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").Value = "427.83"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("D35").Value = "158.60"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "18.98"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "17.97"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("D40").Value = "0.299"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "4.35"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "130.85"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "0.0718"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "0.557"
$ws.Range("D50").Value = "0.0922"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +2.20%  "
